$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "OOTB Domain Groups"
$newSheet.Range("A1").Value = "1and1.com"
$newSheet.Activate()
$newSheet.Range("A23").Select()

$wsWarmup = $wb.Worksheets.Item("Warmup Plan")
$wsWarmup.Activate()
